# Update "paises" (countries) COVID-19 stats sheet with the latest data pull.
#
# Net effect of this update:
#   - The "datos actualizados" timestamp in A1 moves from 00:57 to 02:14.
#   - Several countries' case/death/recovery counters are refreshed with
#     newer figures (rows keyed by their row number below).
#   - Because the sheet is kept sorted by "Casos totales" (col B) descending,
#     Chequia's refreshed total now overtakes Bielorrusia (rows 50/51 swap
#     labels), and Curazao's refreshed total now overtakes Gibraltar
#     (rows 180/181 swap labels) while keeping the rest of each row's data
#     attached to the correct country.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 02:14"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 7600353
$ws.Range("C4").Value = 48432
$ws.Range("D4").Value = 4814551
$ws.Range("E4").Value = 2571531
$ws.Range("G4").Value = 749
$ws.Range("H4").Value = 214271

# --- Row 6: Brasil ---
$ws.Range("E6").Value = 512248
$ws.Range("G6").Value = 580
$ws.Range("H6").Value = 146011

# --- Row 29: Canada ---
$ws.Range("B29").Value = 164471
$ws.Range("C29").Value = 1812
$ws.Range("D29").Value = 138867
$ws.Range("E29").Value = 16142
$ws.Range("G29").Value = 53
$ws.Range("H29").Value = 9462

# --- Row 50: becomes Chequia (refreshed figures overtake Bielorrusia) ---
$ws.Range("A50").Value = "Chequia"
$ws.Range("B50").Value = 80605
$ws.Range("C50").Value = 2554
$ws.Range("D50").Value = 39073
$ws.Range("E50").Value = 40821
$ws.Range("G50").Value = 12
$ws.Range("H50").Value = 711

# --- Row 51: becomes Bielorrusia (keeps its previous figures) ---
$ws.Range("A51").Value = "Bielorrusia"
$ws.Range("B51").Value = 79852
$ws.Range("C51").Value = 431
$ws.Range("D51").Value = 75148
$ws.Range("E51").Value = 3853
$ws.Range("G51").Value = 7
$ws.Range("H51").Value = 851

# --- Row 95: Noruega ---
$ws.Range("B95").Value = 14362
$ws.Range("C95").Value = 78
$ws.Range("E95").Value = 2897

# --- Row 113: Zimbabue ---
$ws.Range("B113").Value = 7885
$ws.Range("C113").Value = 27
$ws.Range("D113").Value = 6327
$ws.Range("E113").Value = 1330

# --- Row 115: Mauritania ---
$ws.Range("B115").Value = 7517
$ws.Range("C115").Value = 6
$ws.Range("D115").Value = 7174
$ws.Range("E115").Value = 181

# --- Row 130: Surinam ---
$ws.Range("B130").Value = 4924
$ws.Range("C130").Value = 25
$ws.Range("D130").Value = 4723
$ws.Range("E130").Value = 95
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 106

# --- Row 167: Niger ---
$ws.Range("B167").Value = 1200
$ws.Range("C167").Value = 2
$ws.Range("E167").Value = 16

# --- Row 169: Santo Tome y Principe ---
$ws.Range("B169").Value = 913
$ws.Range("C169").Value = 2
$ws.Range("E169").Value = 12

# --- Row 180: becomes Curazao (refreshed figures overtake Gibraltar) ---
$ws.Range("A180").Value = "Curazao"
$ws.Range("B180").Value = 429
$ws.Range("C180").Value = 18
$ws.Range("D180").Value = 213
$ws.Range("E180").Value = 215
$ws.Range("H180").Value = 1

# --- Row 181: becomes Gibraltar (keeps its previous figures) ---
$ws.Range("A181").Value = "Gibraltar"
$ws.Range("B181").Value = 428
$ws.Range("C181").Value = 12
$ws.Range("D181").Value = 358
$ws.Range("E181").Value = 70
$ws.Range("H181").Value = 0
